$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Antibodies")
$ws.Activate()

$ws.Range("A4").Value = "Acme mAb 1"

$ws.Range("B5").Value = "Homo sapiens"
$ws.Range("B6").Value = ""
$ws.Range("B7").Value = "Mu musculus"
$ws.Range("B8").Value = "Coronavirus"

$ws.Range("C9").Value = "IgA2"
$ws.Range("C10").Value = ""
$ws.Range("C11").Value = "Ig"
